# Correct / add file extensions for puzzle image links in column F (URL).
# Puzzle_109 and Puzzle_112-118 were actually saved as .JPG (uppercase),
# and Puzzle_129 / Puzzle_147 turned out to be .jpeg rather than .jpg.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    110 = "static/images/Puzzle_109.JPG"
    113 = "static/images/Puzzle_112.JPG"
    114 = "static/images/Puzzle_113.JPG"
    115 = "static/images/Puzzle_114.JPG"
    116 = "static/images/Puzzle_115.JPG"
    117 = "static/images/Puzzle_116.JPG"
    118 = "static/images/Puzzle_117.JPG"
    119 = "static/images/Puzzle_118.JPG"
    130 = "static/images/Puzzle_129.jpeg"
    148 = "static/images/Puzzle_147.jpeg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}

# Restore the view state: selection on the H column (computed JSON) for the
# frozen lower pane, scrolled down near the bottom of the data.
$ws.Range("H2:H176").Select()
